# Pricing model update: the client purchased the "Basic Datamart" add-on,
# so flip the corresponding Inputs flag from "N" to "Y". All downstream
# formulas (Calcs, Outputs_External, Outputs_Internal, Outputs_Timeline)
# recalculate automatically from this single input change.

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Range("D8").Value = "Y"
$wsInputs.Range("D9").Select()

$wsTimeline = $wb.Worksheets.Item("Outputs_Timeline")
$wsTimeline.Activate()
$wsTimeline.Range("O9:O44").Select()
